# Update "Todesfälle & Sterblichkeit" table on sheet "Todesfälle und Fallsterblichkei"
# from KW 47 to KW 48 (destatis update), per commit "update sterbedaten kw 48 destatis".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Todesfälle und Fallsterblichkei")

# Header: previous current-week column becomes header for new current week
$ws.Range("C1").Value = "KW 48"

# Todesfälle & Sterblichkeit - Altersgruppen
# 0 bis 59 Jahre
$ws.Range("B2").Value = "83 ( 0,1%)"
$ws.Range("C2").Value = "82 ( 0,1%)"
$ws.Range("D2").Value = " -1,2%"

# 60 bis 79 Jahre
$ws.Range("B3").Value = "688 ( 3,8%)"
$ws.Range("C3").Value = "688 ( 3,9%)"
$ws.Range("D3").Value = "  0,0%"

# 80 Jahre +
$ws.Range("B4").Value = "1876 (17,4%)"
$ws.Range("C4").Value = "1985 (16,0%)"
$ws.Range("D4").Value = "  5,8%"

# Gesamt
$ws.Range("B5").Value = "2648 ( 2,1%)"
$ws.Range("C5").Value = "2757 ( 2,2%)"
$ws.Range("D5").Value = "  4,1%"

# Übersterblichkeit - Altersgruppen
# 0 bis 59 Jahre
$ws.Range("B7").Value = "-32 (-1,9%)"
$ws.Range("C7").Value = "-102 (-6,0%)"
$ws.Range("D7").Value = "218,8%"

# 60 bis 79 Jahre
$ws.Range("B8").Value = "126 ( 2,1%)"
$ws.Range("C8").Value = "117 ( 1,9%)"
$ws.Range("D8").Value = " -7,1%"

# 80 Jahre +
$ws.Range("B9").Value = "1753 (17,2%)"
$ws.Range("C9").Value = "2510 (24,5%)"
$ws.Range("D9").Value = " 43,2%"

# Gesamt
$ws.Range("B10").Value = "1846 (10,3%)"
$ws.Range("C10").Value = "2525 (13,9%)"
$ws.Range("D10").Value = " 36,8%"
